$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-01 18:24:11"

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
